# Regenerate orders with updated distance/size labels.
#
# The experiment's distance conditions were renamed (D51->D55, D64->D69,
# D80->D86) and the larger "Size" level was renamed (S30->S31). These
# tokens appear throughout the order sheet: in the Condition labels
# (e.g. "Face01_D51_S25"), the left/right stimulus filenames
# (e.g. "Face01_D51_S25_l.png"), and the standalone Distance/Size columns
# (e.g. "D51", "S30"). Every occurrence is the same literal substring, so
# a straightforward find/replace over the whole used range reproduces the
# rename consistently without touching the numeric/boolean columns
# (Trial, Duration_Seconds, Is_Repeat, Block, ConditionID) or the Face
# column, none of which contain these tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$used.Replace("D51", "D55")
$used.Replace("D64", "D69")
$used.Replace("D80", "D86")
$used.Replace("S30", "S31")
